$wb = $excel.ActiveWorkbook

# The "Remove Incomplete Records" worksheet holds the test-run data that
# needs to be filled in (rows 3-6, columns C:L) up to the 20% split.
$ws = $wb.Worksheets.Item("Remove Incomplete Records")

$data = @{
    3 = @(87.18, 79.49, 82.05, 79.49, 84.62, 69.23, 76.92, 74.36, 82.05, 76.92)
    4 = @(74.03, 81.82, 76.62, 74.03, 76.62, 68.83, 70.13, 88.31, 77.92, 77.92)
    5 = @(72.41, 81.90, 81.90, 77.59, 81.03, 77.59, 78.45, 77.59, 74.14, 80.17)
    6 = @(74.68, 74.68, 75.97, 76.62, 74.03, 78.57, 74.03, 74.03, 81.17, 76.62)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i   # column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Recalculate so MIN/MAX/AVERAGE formulas in columns M:O pick up the new
# values (and so Excel drops the stale ignoredErrors for O3:O21 /
# A3:A21,C2:L2 once the #DIV/0! errors are gone).
$excel.CalculateFullRebuild()

# Make "Remove Incomplete Records" the active sheet/tab and set its
# selection, mirroring the saved view state.
$ws.Activate()
$ws.Range("C7").Select()
